$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: "TICKER " -> "Ticker ", "NAME" -> "Name"
$ws.Range("A1").Value = "Ticker "
$ws.Range("B1").Value = "Name"

# Delete the "Codename" column (C) and the "Inception Dates" column (D)
$ws.Columns.Item(4).Delete()
$ws.Columns.Item(3).Delete()

# Move selection to match target state
$ws.Range("B29").Select()
